{"js": "const pairs = [\n  [\"68-28=40\", \"14+23=37\"],\n  [\"34+43=77\", \"15+7=22\"],\n  [\"64-41=23\", \"71-45=26\"],\n  [\"38+1=39\", \"10+74=84\"],\n  [\"76-74=2\", \"90-13=77\"],\n  [\"2+25=27\", \"0+28=28\"],\n  [\"31+1=32\", \"8+75=83\"],\n  [\"99-61=38\", \"57-48=9\"],\n  [\"14+51=65\", \"95+1=96\"],\n  [\"43-42=1\", \"2+35=37\"],\n  [\"78+14=92\", \"88-14=74\"],\n  [\"52+20=72\", \"85+12=97\"],\n  [\"25+51=76\", \"56+11=67\"],\n  [\"97-48=49\", \"24-7=17\"],\n  [\"31+65=96\", \"22+57=79\"],\n  [\"84-59=25\", \"56-15=41\"],\n  [\"53+32=85\", \"83-34=49\"],\n  [\"70-21=49\", \"66+14=80\"],\n  [\"8+52=60\", \"88-67=21\"],\n  [\"60-3=57\", \"77-8=69\"],\n  [\"11+46=57\", \"55-15=40\"],\n  [\"53-48=5\", \"75+13=88\"],\n  [\"57-10=47\", \"99-84=15\"],\n  [\"17+14=31\", \"73-15=58\"],\n  [\"54-26=28\", \"55-4=51\"],\n  [\"1+33=34\", \"39-38=1\"],\n  [\"77-6=71\", \"8+39=47\"],\n  [\"41+47=88\", \"47-37=10\"],\n  [\"85-85=0\", \"37-14=23\"],\n  [\"42+40=82\", \"90-87=3\"],\n  [\"55+21=76\", \"29-8=21\"],\n  [\"52+23=75\", \"51+43=94\"],\n  [\"77-41=36\", \"96-41=55\"],\n  [\"37+52=89\", \"34-33=1\"],\n  [\"33-12=21\", \"43+39=82\"],\n  [\"7+59=66\", \"20+38=58\"],\n  [\"28+20=48\", \"57-37=20\"],\n  [\"17+19=36\", \"67-61=6\"],\n  [\"32+63=95\", \"19+54=73\"],\n  [\"32+67=99\", \"4+45=49\"],\n  [\"76-75=1\", \"32-29=3\"],\n  [\"12-11=1\", \"33+3=36\"],\n  [\"60-12=48\", \"71-55=16\"],\n  [\"21+4=25\", \"62-48=14\"],\n  [\"36+35=71\", \"97-80=17\"],\n  [\"83-44=39\", \"61-23=38\"],\n  [\"60-8=52\", \"48-11=37\"],\n  [\"52+7=59\", \"29+58=87\"],\n  [\"92+7=99\", \"87-72=15\"],\n  [\"58+37=95\", \"41+33=74\"],\n  [\"27-8=19\", \"24+7=31\"],\n  [\"55-47=8\", \"48+13=61\"],\n  [\"46+22=68\", \"80-25=55\"],\n  [\"63+24=87\", \"52-21=31\"],\n  [\"59-43=16\", \"60-2=58\"],\n  [\"54+18=72\", \"30-17=13\"],\n  [\"23+39=62\", \"12+70=82\"],\n  [\"50+41=91\", \"58-44=14\"],\n  [\"25+24=49\", \"91-71=20\"],\n  [\"4+60=64\", \"83-68=15\"],\n  [\"84-46=38\", \"90-10=80\"],\n  [\"20+52=72\", \"10-5=5\"],\n  [\"7+28=35\", \"66-52=14\"],\n  [\"55-45=10\", \"46-17=29\"],\n  [\"63-11=52\", \"31+3=34\"],\n  [\"46-28=18\", \"76-47=29\"],\n  [\"40+30=70\", \"90-53=37\"],\n  [\"53+21=74\", \"67-45=22\"],\n  [\"61-21=40\", \"19+3=22\"],\n  [\"24+49=73\", \"5+60=65\"],\n  [\"63-63=0\", \"35-33=2\"],\n  [\"91-21=70\", \"15-3=12\"],\n  [\"5+62=67\", \"18+46=64\"],\n  [\"53+15=68\", \"15+55=70\"],\n  [\"10+46=56\", \"27+18=45\"],\n  [\"80-27=53\", \"6-3=3\"],\n  [\"45+3=48\", \"75-29=46\"],\n  [\"83-23=60\", \"62+37=99\"],\n  [\"39+52=91\", \"56+8=64\"],\n  [\"82-35=47\", \"85+11=96\"],\n  [\"85-23=62\", \"35-16=19\"],\n  [\"37+60=97\", \"27-6=21\"],\n  [\"30-26=4\", \"1+70=71\"],\n  [\"25+38=63\", \"50+2=52\"],\n  [\"96-3=93\", \"39-1=38\"],\n  [\"51+28=79\", \"98-65=33\"],\n  [\"44+36=80\", \"67+30=97\"],\n  [\"6+58=64\", \"81-14=67\"],\n  [\"16-6=10\", \"79+5=84\"],\n  [\"62-34=28\", \"76-13=63\"],\n  [\"30-12=18\", \"31-29=2\"],\n  [\"72-22=50\", \"74+7=81\"],\n  [\"77-46=31\", \"75-1=74\"],\n  [\"63-9=54\", \"98-54=44\"],\n  [\"75+23=98\", \"92-63=29\"],\n  [\"17-15=2\", \"7+64=71\"],\n  [\"7+70=77\", \"31+56=87\"],\n  [\"76+3=79\", \"17+15=32\"],\n  [\"45+51=96\", \"38+59=97\"],\n  [\"4+78=82\", \"2+11=13\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"68-28=40\", \"14+23=37\"),\n    @(\"34+43=77\", \"15+7=22\"),\n    @(\"64-41=23\", \"71-45=26\"),\n    @(\"38+1=39\", \"10+74=84\"),\n    @(\"76-74=2\", \"90-13=77\"),\n    @(\"2+25=27\", \"0+28=28\"),\n    @(\"31+1=32\", \"8+75=83\"),\n    @(\"99-61=38\", \"57-48=9\"),\n    @(\"14+51=65\", \"95+1=96\"),\n    @(\"43-42=1\", \"2+35=37\"),\n    @(\"78+14=92\", \"88-14=74\"),\n    @(\"52+20=72\", \"85+12=97\"),\n    @(\"25+51=76\", \"56+11=67\"),\n    @(\"97-48=49\", \"24-7=17\"),\n    @(\"31+65=96\", \"22+57=79\"),\n    @(\"84-59=25\", \"56-15=41\"),\n    @(\"53+32=85\", \"83-34=49\"),\n    @(\"70-21=49\", \"66+14=80\"),\n    @(\"8+52=60\", \"88-67=21\"),\n    @(\"60-3=57\", \"77-8=69\"),\n    @(\"11+46=57\", \"55-15=40\"),\n    @(\"53-48=5\", \"75+13=88\"),\n    @(\"57-10=47\", \"99-84=15\"),\n    @(\"17+14=31\", \"73-15=58\"),\n    @(\"54-26=28\", \"55-4=51\"),\n    @(\"1+33=34\", \"39-38=1\"),\n    @(\"77-6=71\", \"8+39=47\"),\n    @(\"41+47=88\", \"47-37=10\"),\n    @(\"85-85=0\", \"37-14=23\"),\n    @(\"42+40=82\", \"90-87=3\"),\n    @(\"55+21=76\", \"29-8=21\"),\n    @(\"52+23=75\", \"51+43=94\"),\n    @(\"77-41=36\", \"96-41=55\"),\n    @(\"37+52=89\", \"34-33=1\"),\n    @(\"33-12=21\", \"43+39=82\"),\n    @(\"7+59=66\", \"20+38=58\"),\n    @(\"28+20=48\", \"57-37=20\"),\n    @(\"17+19=36\", \"67-61=6\"),\n    @(\"32+63=95\", \"19+54=73\"),\n    @(\"32+67=99\", \"4+45=49\"),\n    @(\"76-75=1\", \"32-29=3\"),\n    @(\"12-11=1\", \"33+3=36\"),\n    @(\"60-12=48\", \"71-55=16\"),\n    @(\"21+4=25\", \"62-48=14\"),\n    @(\"36+35=71\", \"97-80=17\"),\n    @(\"83-44=39\", \"61-23=38\"),\n    @(\"60-8=52\", \"48-11=37\"),\n    @(\"52+7=59\", \"29+58=87\"),\n    @(\"92+7=99\", \"87-72=15\"),\n    @(\"58+37=95\", \"41+33=74\"),\n    @(\"27-8=19\", \"24+7=31\"),\n    @(\"55-47=8\", \"48+13=61\"),\n    @(\"46+22=68\", \"80-25=55\"),\n    @(\"63+24=87\", \"52-21=31\"),\n    @(\"59-43=16\", \"60-2=58\"),\n    @(\"54+18=72\", \"30-17=13\"),\n    @(\"23+39=62\", \"12+70=82\"),\n    @(\"50+41=91\", \"58-44=14\"),\n    @(\"25+24=49\", \"91-71=20\"),\n    @(\"4+60=64\", \"83-68=15\"),\n    @(\"84-46=38\", \"90-10=80\"),\n    @(\"20+52=72\", \"10-5=5\"),\n    @(\"7+28=35\", \"66-52=14\"),\n    @(\"55-45=10\", \"46-17=29\"),\n    @(\"63-11=52\", \"31+3=34\"),\n    @(\"46-28=18\", \"76-47=29\"),\n    @(\"40+30=70\", \"90-53=37\"),\n    @(\"53+21=74\", \"67-45=22\"),\n    @(\"61-21=40\", \"19+3=22\"),\n    @(\"24+49=73\", \"5+60=65\"),\n    @(\"63-63=0\", \"35-33=2\"),\n    @(\"91-21=70\", \"15-3=12\"),\n    @(\"5+62=67\", \"18+46=64\"),\n    @(\"53+15=68\", \"15+55=70\"),\n    @(\"10+46=56\", \"27+18=45\"),\n    @(\"80-27=53\", \"6-3=3\"),\n    @(\"45+3=48\", \"75-29=46\"),\n    @(\"83-23=60\", \"62+37=99\"),\n    @(\"39+52=91\", \"56+8=64\"),\n    @(\"82-35=47\", \"85+11=96\"),\n    @(\"85-23=62\", \"35-16=19\"),\n    @(\"37+60=97\", \"27-6=21\"),\n    @(\"30-26=4\", \"1+70=71\"),\n    @(\"25+38=63\", \"50+2=52\"),\n    @(\"96-3=93\", \"39-1=38\"),\n    @(\"51+28=79\", \"98-65=33\"),\n    @(\"44+36=80\", \"67+30=97\"),\n    @(\"6+58=64\", \"81-14=67\"),\n    @(\"16-6=10\", \"79+5=84\"),\n    @(\"62-34=28\", \"76-13=63\"),\n    @(\"30-12=18\", \"31-29=2\"),\n    @(\"72-22=50\", \"74+7=81\"),\n    @(\"77-46=31\", \"75-1=74\"),\n    @(\"63-9=54\", \"98-54=44\"),\n    @(\"75+23=98\", \"92-63=29\"),\n    @(\"17-15=2\", \"7+64=71\"),\n    @(\"7+70=77\", \"31+56=87\"),\n    @(\"76+3=79\", \"17+15=32\"),\n    @(\"45+51=96\", \"38+59=97\"),\n    @(\"4+78=82\", \"2+11=13\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $result = $find.Execute([ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, \"wdReplaceAll\")\n    if (-not $result) {\n        throw \"Replace failed for: $old\"\n    }\n}"}
